$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps being treated as text (it already stores values
# like "1.000" / "29.873.95" as text, not numbers), so updated values do not get
# silently reinterpreted as numbers/dates by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.873.95"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "1.893.67"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "0.7822"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("D6").Value = "243.90"
$ws.Range("E6").Value = "  +1.03%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "0.3134"
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("D10").Value = "0.07336"
$ws.Range("D11").Value = "0.08096"
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").Value = "0.7732"
$ws.Range("E12").Value = "  +1.72%  "
$ws.Range("D13").Value = "5.473"
$ws.Range("E13").Value = "  +4.16%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.907.44"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "93.78"
$ws.Range("E15").Value = "  +2.06%  "
$ws.Range("D16").Value = "6.221"
$ws.Range("E16").Value = "  +5.73%  "
$ws.Range("D17").Value = "29.885.65"
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("D18").Value = "13.95"
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("D19").Value = "246.96"
$ws.Range("E19").Value = "  +2.13%  "
$ws.Range("D20").Value = "0.000007822"
$ws.Range("E20").Value = "  +2.05%  "
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").Value = "8.114"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "2.137.70"
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").Value = "0.1594"
$ws.Range("E25").Value = "  -2.12%  "
$ws.Range("D26").Value = "9.428"
$ws.Range("E26").Value = "  +1.81%  "
$ws.Range("D27").Value = "164.09"
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("E29").Value = "  -0.80%  "
$ws.Range("E30").Value = "  +2.52%  "
$ws.Range("D31").Value = "1.542"
$ws.Range("E31").Value = "  +0.62%  "
$ws.Range("D32").Value = "4.479"
$ws.Range("E32").Value = "  +2.70%  "
$ws.Range("D33").Value = "0.05554"
$ws.Range("E33").Value = "  -1.49%  "
$ws.Range("D34").Value = "4.056"
$ws.Range("E34").Value = "  +0.97%  "
$ws.Range("D35").Value = "1.240"
$ws.Range("E35").Value = "  -1.15%  "
$ws.Range("D36").Value = "0.7517"
$ws.Range("E36").Value = "  +2.83%  "
$ws.Range("D37").Value = "0.9994"
$ws.Range("E37").Value = "  -0.36%  "
$ws.Range("E38").Value = "  +1.81%  "
$ws.Range("D39").Value = "0.01938"
$ws.Range("E39").Value = "  +2.15%  "
$ws.Range("D40").Value = "2.797"
$ws.Range("E40").Value = "  +0.91%  "
$ws.Range("D41").Value = "1.138.91"
$ws.Range("E41").Value = "  +12.11%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "74.26"
$ws.Range("E42").Value = "  +3.17%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "0.4460"
$ws.Range("E43").Value = "  +1.97%  "
$ws.Range("D44").Value = "5.968"
$ws.Range("D45").Value = "0.8525"
$ws.Range("E45").Value = "  +2.03%  "
$ws.Range("D46").Value = "0.9999"
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("D47").Value = "1.889"
$ws.Range("E47").Value = "  +2.23%  "
$ws.Range("D48").Value = "102.26"
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("D49").Value = "3.060"
$ws.Range("E49").Value = "  +6.03%  "
$ws.Range("D50").Value = "7.513"
$ws.Range("E50").Value = "  +2.06%  "
$ws.Range("D51").Value = "9.752"
$ws.Range("E51").Value = "  -0.45%  "
